$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 79 new "hội chiếu tại Phu Thê" conjunction phrases added at the end of
# the sheet (rows 4224-4302), mirrored identically into columns A and B.
$values = @(
    'Sát Phá Tham hội chiếu tại Phu Thê',
    'Tử Phủ Vũ Tướng hội chiếu tại Phu Thê',
    'Cơ Nguyệt Đồng Lương hội chiếu tại Phu Thê',
    'Cự Nhật hội chiếu tại Phu Thê',
    'Kình Đà hội chiếu tại Phu Thê',
    'Xương Khúc hội chiếu tại Phu Thê',
    'Hoả Linh hội chiếu tại Phu Thê',
    'Không Kiếp hội chiếu tại Phu Thê',
    'Quang Quý hội chiếu tại Phu Thê',
    'Tả Hữu hội chiếu tại Phu Thê',
    'Song Hao hội chiếu tại Phu Thê',
    'Tang Hổ hội chiếu tại Phu Thê',
    'Khốc Hư hội chiếu tại Phu Thê',
    'Hình Riêu hội chiếu tại Phu Thê',
    'Thai Toạ hội chiếu tại Phu Thê',
    'Đào Hồng hội chiếu tại Phu Thê',
    'Ấn Phù hội chiếu tại Phu Thê',
    'Song Hao Quyền Lộc Kiếp Hoả hội chiếu tại Phu Thê',
    'Tử Phủ Vũ Tướng Xương Khúc Khôi Việt Tả Hữu Khoa Quyền Lộc Long hội chiếu tại Phu Thê',
    'Tử Phủ Vũ Tướng Tả Hữu Khoa Quyền Lộc Long Phượng hội chiếu tại Phu Thê',
    'Tử Khúc Phá Dương Đà hội chiếu tại Phu Thê',
    'Cơ Nguyệt Đồng Lương gặp Xương Khúc Tả Hữu hội chiếu tại Phu Thê',
    'Cơ Nguyệt Đồng Lương Khoa Tả Hữu Quang Quý Quan Phúc hội chiếu tại Phu Thê',
    'Sát Quyền hội chiếu tại Phu Thê',
    'Lộc Mã hội chiếu tại Phu Thê',
    'Kiếp Hư Hao Quyền hội chiếu tại Phu Thê',
    'Tuế Hổ Phù Xương Khúc hội chiếu tại Phu Thê',
    'Xương Khúc Tấu Long Phượng hội chiếu tại Phu Thê',
    'Đào Hồng Riêu Tấu Cơ Vũ hội chiếu tại Phu Thê',
    'Binh Hình Tướng Ấn hội chiếu tại Phu Thê',
    'Hổ Tấu hội chiếu tại Phu Thê',
    'Hình Riêu Y hội chiếu tại Phu Thê',
    'Mã Hỏa Linh hội chiếu tại Phu Thê',
    'Thai Tọa Hồng Đào hội chiếu tại Phu Thê',
    'Tả Hữu Không Kiếp hội chiếu tại Phu Thê',
    'Tả Hữu Binh Tướng hội chiếu tại Phu Thê',
    'Đào Quyền hội chiếu tại Phu Thê',
    'Đào Hồng Tả Cái Hữu hội chiếu tại Phu Thê',
    'Quan Phúc Quang Tấu hội chiếu tại Phu Thê',
    'Đào Hồng Xương Khúc hội chiếu tại Phu Thê',
    'Đào Tử Phủ hội chiếu tại Phu Thê',
    'Hổ Kình Sát hội chiếu tại Phu Thê',
    'Hổ Tang Kiếp hội chiếu tại Phu Thê',
    'Hổ Tang Không Kiếp hội chiếu tại Phu Thê',
    'Âm Dương Lương hội chiếu tại Phu Thê',
    'Cơ Lương Gia Hội hội chiếu tại Phu Thê',
    'Nhật Chiếu Lôi Môn hội chiếu tại Phu Thê',
    'Tả Hữu Xương Khúc hội chiếu tại Phu Thê',
    'Tham Linh Triều Viên hội chiếu tại Phu Thê',
    'Cự Hỏa Linh hội chiếu tại Phu Thê',
    'Hồng Đào Kỵ hội chiếu tại Phu Thê',
    'Cự Đồng Hình hội chiếu tại Phu Thê',
    'Lương Phá hội chiếu tại Phu Thê',
    'Khôi Việt hội chiếu tại Phu Thê',
    'Tham Vũ Hỏa hội chiếu tại Phu Thê',
    'Cự Kỵ hội chiếu tại Phu Thê',
    'Kình Đà Hỏa hội chiếu tại Phu Thê',
    'Kình Đà Không Kiếp hội chiếu tại Phu Thê',
    'Đào Hồng Hỷ hội chiếu tại Phu Thê',
    'Đào Hồng Kiếp Sát hội chiếu tại Phu Thê',
    'Đào Hồng Kỵ hội chiếu tại Phu Thê',
    'Đào Hồng Tả Phù Hữu Bật hội chiếu tại Phu Thê',
    'Đào Hồng Tả Phù Hữu Bật Khoa Quyền Lộc hội chiếu tại Phu Thê',
    'Đào Hồng Tả Phù Hữu Bật Khoa Quyền Lộc Long Trì Phượng Các hội chiếu tại Phu Thê',
    'Lưu Hà Kiếp Sát hội chiếu tại Phu Thê',
    'Phục Không Kiếp hội chiếu tại Phu Thê',
    'Xương Khúc Khôi Việt hội chiếu tại Phu Thê',
    'Khoa Quyền hội chiếu tại Phu Thê',
    'Lộc Quyền hội chiếu tại Phu Thê',
    'Tướng Binh Đào Hồng hội chiếu tại Phu Thê',
    'Lương Khốc Tuế hội chiếu tại Phu Thê',
    'Tướng Binh hội chiếu tại Phu Thê',
    'Cự Tang hội chiếu tại Phu Thê',
    'Cự Hỏa hội chiếu tại Phu Thê',
    'Phá Hình Kỵ hội chiếu tại Phu Thê',
    'Tang Trực Tuế hội chiếu tại Phu Thê',
    'Tuế Xương Khúc hội chiếu tại Phu Thê',
    'Tử Vi Tả Hữu hội chiếu tại Phu Thê',
    'Tử Vi Tang Tả Hữu hội chiếu tại Phu Thê'
)

$startRow = 4224
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$lastRow = $startRow + $values.Count - 1

# Match the author's final view state: scrolled down with the newly added
# block (B4224:B4302) selected.
$win = $excel.ActiveWindow
$ws.Range("B" + $startRow + ":B" + $lastRow).Select() | Out-Null
$win.ScrollRow = 4216
